$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "0FAP5B"
$ws.Range("B8").Value = "Almohadilla + Chip Epson"
$ws.Range("C8").Value = "SC 23MB F100 F130 F160 F170"
$ws.Range("D8").Value = 35000
$ws.Range("E8").Value = 200000
$ws.Range("F8").Value = 11
$ws.Range("G8").Value = 16
$ws.Range("H8").Formula = "=(E8-D8)*G8"
$ws.Range("I8").Formula = "=D8*F8"
$ws.Range("J8").Value = 385000
